$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11609
$ws.Range("J70").Value = 12787.777
$ws.Range("L70").Value = 38363.331
$ws.Range("N70").Value = -38903.331

$ws.Range("H73").Value = 11609
$ws.Range("J73").Value = 12787.777
$ws.Range("L73").Value = 38363.331
$ws.Range("N73").Value = -40235.331

$ws.Range("H132").Value = 554.3012
$ws.Range("I132").Value = 500.21518
$ws.Range("K132").Value = 1500.64554
$ws.Range("M132").Value = 1029.35446

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2907620
$ws.Range("I2").Value = 4651692
$ws.Range("K2").Value = 4651692
$ws.Range("M2").Value = -4651579

$ws.Range("H32").Value = 3408.0322
$ws.Range("I32").Value = 2442.6853
$ws.Range("K32").Value = 2442.6853
$ws.Range("M32").Value = -2155.6853

$ws.Range("H46").Value = 11607.2
$ws.Range("J46").Value = 12866.333
$ws.Range("L46").Value = 12866.333
$ws.Range("N46").Value = -13504.333

$ws.Range("H63").Value = 3199.6
$ws.Range("I63").Value = 9998
$ws.Range("K63").Value = 9998
$ws.Range("M63").Value = -9312

$ws.Range("H66").Value = 3199.6
$ws.Range("I66").Value = 9998
$ws.Range("K66").Value = 49990
$ws.Range("M66").Value = -46558

$ws.Range("H101").Value = 16773
$ws.Range("J101").Value = 16773
$ws.Range("L101").Value = 16773
$ws.Range("N101").Value = -23263

$ws.Range("H116").Value = 2907620
$ws.Range("I116").Value = 4651692
$ws.Range("K116").Value = 4651692
$ws.Range("M116").Value = -4649398

$ws.Range("H123").Value = 72000
$ws.Range("J123").Value = 72000
$ws.Range("L123").Value = 72000
$ws.Range("N123").Value = -81800

$ws.Range("H125").Value = 49997.5
$ws.Range("J125").Value = 49997.5
$ws.Range("L125").Value = 49997.5
$ws.Range("N125").Value = -59837.5

$ws.Range("H135").Value = 23957
$ws.Range("J135").Value = 23957
$ws.Range("L135").Value = 23957
$ws.Range("N135").Value = -34097

$ws.Range("H139").Value = 52571.668
$ws.Range("J139").Value = 52571.668
$ws.Range("L139").Value = 52571.668
$ws.Range("N139").Value = -62851.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2907620
$ws.Range("I3").Value = 4651692
$ws.Range("K3").Value = 4651692
$ws.Range("M3").Value = -4651578

$ws.Range("H81").Value = 19396
$ws.Range("J81").Value = 19396
$ws.Range("L81").Value = 19396
$ws.Range("N81").Value = -21518

$ws.Range("H84").Value = 19396
$ws.Range("J84").Value = 19396
$ws.Range("L84").Value = 58188
$ws.Range("N84").Value = -68796

$ws.Range("H107").Value = 1648.4546
$ws.Range("I107").Value = 1413.2
$ws.Range("K107").Value = 1413.2
$ws.Range("M107").Value = 506.8

$ws.Range("H110").Value = 49999
$ws.Range("J110").Value = 49999
$ws.Range("L110").Value = 49999
$ws.Range("N110").Value = -58179

$ws.Range("H134").Value = 8655.317999999999
$ws.Range("I134").Value = 13435
$ws.Range("K134").Value = 40305
$ws.Range("M134").Value = -37770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3568.8333
$ws.Range("I31").Value = 3756.889
$ws.Range("K31").Value = 3756.889
$ws.Range("M31").Value = -3461.889

$ws.Range("H34").Value = 3568.8333
$ws.Range("I34").Value = 3756.889
$ws.Range("K34").Value = 3756.889
$ws.Range("M34").Value = -3554.889

$ws.Range("H43").Value = 22500
$ws.Range("J43").Value = 22500
$ws.Range("L43").Value = 22500
$ws.Range("N43").Value = -22868

$ws.Range("H86").Value = 71430110
$ws.Range("I86").Value = 90910210
$ws.Range("K86").Value = 90910210
$ws.Range("M86").Value = -90909087

$ws.Range("H89").Value = 71430110
$ws.Range("I89").Value = 90910210
$ws.Range("K89").Value = 454551050
$ws.Range("M89").Value = -454545434

$ws.Range("H94").Value = 894.5
$ws.Range("I94").Value = 818.5
$ws.Range("J94").Value = 932.5
$ws.Range("K94").Value = 818.5
$ws.Range("L94").Value = 932.5
$ws.Range("M94").Value = -367.5
$ws.Range("N94").Value = -1834.5

$ws.Range("H95").Value = 26655.334
$ws.Range("J95").Value = 26655.334
$ws.Range("L95").Value = 26655.334
$ws.Range("N95").Value = -32147.334

$ws.Range("H101").Value = 22500
$ws.Range("J101").Value = 22500
$ws.Range("L101").Value = 22500
$ws.Range("N101").Value = -28990

$ws.Range("H134").Value = 1153.7273
$ws.Range("I134").Value = 1142.9375
$ws.Range("J134").Value = 1499
$ws.Range("K134").Value = 3428.8125
$ws.Range("L134").Value = 4497
$ws.Range("M134").Value = -893.8125
$ws.Range("N134").Value = -9567

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 798.1667

$ws.Range("H60").Value = 1650
$ws.Range("I60").Value = 300
$ws.Range("J60").Value = 3000
$ws.Range("K60").Value = 900
$ws.Range("L60").Value = 9000
$ws.Range("M60").Value = -649
$ws.Range("N60").Value = -9502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5175
$ws.Range("I70").Value = 5850
$ws.Range("K70").Value = 5850
$ws.Range("M70").Value = -5580

$ws.Range("H73").Value = 5175
$ws.Range("I73").Value = 5850
$ws.Range("K73").Value = 5850
$ws.Range("M73").Value = -4914

$ws.Range("H95").Value = 27500
$ws.Range("J95").Value = 27500
$ws.Range("L95").Value = 27500
$ws.Range("N95").Value = -32992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5665.6665
$ws.Range("I40").Value = 2866.6667
$ws.Range("J40").Value = 8464.666999999999
$ws.Range("K40").Value = 2866.6667
$ws.Range("L40").Value = 8464.666999999999
$ws.Range("M40").Value = -2730.6667
$ws.Range("N40").Value = -8736.666999999999

$ws.Range("H100").Value = 1833.3334
$ws.Range("I100").Value = 1833.3334
$ws.Range("K100").Value = 1833.3334
$ws.Range("M100").Value = -1292.3334

$ws.Range("H136").Value = 3549.3684
$ws.Range("I136").Value = 2159.182
$ws.Range("J136").Value = 5460.875
$ws.Range("K136").Value = 6477.545999999999
$ws.Range("L136").Value = 16382.625
$ws.Range("M136").Value = -3927.545999999999
$ws.Range("N136").Value = -21482.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 39221
$ws.Range("J105").Value = 39221
$ws.Range("L105").Value = 39221
$ws.Range("N105").Value = -46209

$ws.Range("H123").Value = 47599.25
$ws.Range("J123").Value = 47599.25
$ws.Range("L123").Value = 47599.25
$ws.Range("N123").Value = -57399.25

$ws.Range("H136").Value = 12921500
$ws.Range("I136").Value = 16836358
$ws.Range("K136").Value = 50509074
$ws.Range("M136").Value = -50506524
